# Apply the changes described by the commit "Update gh-pages to output
# generated at 456a3b4" to the 广州-漫展信息.xlsx workbook.
#
# The workbook has 4 sheets:
#   展览     (sheet1 / rId1) - Exhibitions
#   演出     (sheet2 / rId2) - Performances
#   本地生活 (sheet3 / rId3) - Local life
#   全部类型 (sheet4 / rId4) - All types (aggregate of the above)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: 展览 (Exhibitions)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets("展览")

$ws1.Range("F2").Value  = 283
$ws1.Range("F3").Value  = 0
$ws1.Range("F6").Value  = 478
$ws1.Range("G9").Value  = 0
$ws1.Range("F10").Value = 0
$ws1.Range("F11").Value = 153
$ws1.Range("F13").Value = 680
$ws1.Range("F14").Value = 0
$ws1.Range("F15").Value = 104
$ws1.Range("F16").Value = 0
$ws1.Range("F17").Value = 0
$ws1.Range("F18").Value = 0
$ws1.Range("F20").Value = 83
$ws1.Range("F23").Value = 1021
$ws1.Range("F25").Value = 535
$ws1.Range("F28").Value = 556
$ws1.Range("F29").Value = 33
$ws1.Range("F31").Value = 0
$ws1.Range("F34").Value = 356
$ws1.Range("F35").Value = 170
$ws1.Range("F36").Value = 225
$ws1.Range("F40").Value = 984
$ws1.Range("F42").Value = 71
$ws1.Range("F43").Value = 0

# ---------------------------------------------------------------------------
# Sheet: 演出 (Performances)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets("演出")

$ws2.Range("C2").Value  = "广州·《唤起心中的能量》唯美国风经典影视作品音乐会"
$ws2.Range("F2").Value  = 3
$ws2.Range("F6").Value  = 7
$ws2.Range("F8").Value  = 7
$ws2.Range("F10").Value = 70
$ws2.Range("F12").Value = 6
$ws2.Range("F21").Value = 5

# ---------------------------------------------------------------------------
# Sheet: 本地生活 (Local life)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets("本地生活")

$ws3.Range("F3").Value = 423
$ws3.Range("F4").Value = 275

# ---------------------------------------------------------------------------
# Sheet: 全部类型 (All types)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets("全部类型")

$ws4.Range("F2").Value = 0
$ws4.Range("F3").Value = 0
$ws4.Range("F5").Value = 283
$ws4.Range("F6").Value = 950
$ws4.Range("F7").Value = 94

$ws4.Range("C9").Value = "广州·《唤起心中的能量》唯美国风经典影视作品音乐会"
$ws4.Range("F9").Value = 0

# Row 11 previously held "萤火虫动漫游戏嘉年华 x KKWORLD2024 快看漫画乐园" which
# was removed from the list. Rows 12-15 shift up into rows 11-14, and a brand
# new row (the "浪漫古典II" concert) is introduced at row 15.
#
# NOTE: the B column holds plain text dates like "2024-07-20" (t="inlineStr"
# in the original file, not real dates). Writing such a string straight into
# .Value makes Excel auto-detect it as a date and reformat the cell, so we
# force the cell to Text format first and clear the resulting formatting
# afterwards so no stray style index is left behind.
$ws4.Range("B11").NumberFormat = "@"
$ws4.Range("B11").Value = "2024-07-20"
$ws4.Range("B11").ClearFormats()
$ws4.Range("C11").Value = "广州·冰兔2024线下live「过去和未来」"
$ws4.Range("D11").Value = "恩宁路265号三层四层自编01 MAO Livehouse广州（永庆坊店）"
$ws4.Range("E11").Value = "2024.07.20 20:00-07.20 22:00"
$ws4.Range("F11").Value = 188
$ws4.Range("G11").Value = 198
$ws4.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=87546"
$ws4.Range("I11").Value = "//i2.hdslb.com/bfs/openplatform/202406/2X09PE1a1718611339266.jpeg"

$ws4.Range("B12").Value = "2024-07-20"
$ws4.Range("C12").Value = "广州·跨越二次元ACG神级动漫世界巡回演唱会"
$ws4.Range("D12").Value = "广州市荔湾区十甫路125号(上下九步行街内)2层 广州平安大戏院"
$ws4.Range("E12").Value = "2024.07.20 19:30-07.20 21:10"
$ws4.Range("F12").Value = 0
$ws4.Range("G12").Value = 280
$ws4.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=85353"
$ws4.Range("I12").Value = "//i1.hdslb.com/bfs/openplatform/202405/4gACWbPh1715223804704.jpeg"

$ws4.Range("B13").NumberFormat = "@"
$ws4.Range("B13").Value = "2024-07-21"
$ws4.Range("B13").ClearFormats()
$ws4.Range("C13").Value = "广州·昨日重现——唯美英文经典歌曲演唱会"
$ws4.Range("D13").Value = "东风中路299号 广州中山纪念堂"
$ws4.Range("E13").Value = "2024.07.21 19:30-07.21 21:30"
$ws4.Range("F13").Value = 0
$ws4.Range("G13").Value = 100
$ws4.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=86802"
$ws4.Range("I13").Value = "//i1.hdslb.com/bfs/openplatform/202405/DR8AvmXe1716802703006.jpeg"

$ws4.Range("B14").Value = "2024-07-21"
$ws4.Range("C14").Value = "广州·燃动!!高梨康治SUMMER LIVE TOUR IN CHINA 2024"
$ws4.Range("D14").Value = "海珠同创汇东一街11号（上冲南约11-2） 声音共和Livehouse"
$ws4.Range("E14").Value = "2024.07.21 14:30-07.21 16:00"
$ws4.Range("F14").Value = 0
$ws4.Range("G14").Value = 280
$ws4.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=87034"
$ws4.Range("I14").Value = "//i1.hdslb.com/bfs/openplatform/202406/LINsP2ui1717741701901.png"

$ws4.Range("B15").NumberFormat = "@"
$ws4.Range("B15").Value = "2024-07-26"
$ws4.Range("B15").ClearFormats()
$ws4.Range("C15").Value = "广州·【早鸟8折】“浪漫古典Ⅱ”百年经典传世名曲烛光音乐会 "
$ws4.Range("D15").Value = "广州市二沙岛晴波路33号  星海音乐厅（交响乐演奏厅）"
$ws4.Range("E15").Value = "2024.07.26 20:00-07.26 21:30"
$ws4.Range("G15").Value = 144
$ws4.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=87726"
$ws4.Range("I15").Value = "//i0.hdslb.com/bfs/openplatform/202406/A8vhVlhn1717575084179.png"

$ws4.Range("F16").Value = 0
$ws4.Range("F17").Value = 153
$ws4.Range("F18").Value = 0
$ws4.Range("F20").Value = 680
$ws4.Range("F21").Value = 542
$ws4.Range("F23").Value = 83
$ws4.Range("F24").Value = 104
$ws4.Range("F25").Value = 194
$ws4.Range("F26").Value = 637
$ws4.Range("F30").Value = 543
$ws4.Range("F32").Value = 1021
$ws4.Range("F33").Value = 0
$ws4.Range("F35").Value = 0
$ws4.Range("F36").Value = 534
$ws4.Range("F38").Value = 33
$ws4.Range("F40").Value = 0
$ws4.Range("F41").Value = 121
$ws4.Range("F42").Value = 0
$ws4.Range("F43").Value = 356
$ws4.Range("F44").Value = 225
$ws4.Range("F45").Value = 192
$ws4.Range("F47").Value = 0
$ws4.Range("F48").Value = 71
$ws4.Range("F50").Value = 0
